$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

$ws.Range("A3").Value = "Testing2"
$ws.Range("A4").Value = "Testing3"
$ws.Range("A2").Value = "Testing4"

$ws.Range("A2").Select()
